$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting the existing rows 83-87 down to 84-88
$ws.Rows.Item(83).Insert()

# Populate the new row 83 with the new weekly record
$ws.Cells.Item(83, 1).Value = 10
$ws.Cells.Item(83, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value = "La Araucanía"
$ws.Cells.Item(83, 4).Value = 45075
$ws.Cells.Item(83, 5).Value = 9
$ws.Cells.Item(83, 6).Value = 100112026
$ws.Cells.Item(83, 7).Value = "Haba"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 75
$ws.Cells.Item(83, 11).Value = 15000
$ws.Cells.Item(83, 12).Value = 15000
$ws.Cells.Item(83, 13).Value = 15000
$ws.Cells.Item(83, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(83, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(83, 16).Value = 600
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
